$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Replace the "Limitations include" bullet (Q2 stub) with:
#    - 4 new sub-bullets answering Q1 (Kickstarter conclusions) at ilvl=1
#    - the relocated _GoBack bookmark paragraph (ListParagraph, no numPr)
#    - the reworded "Limitations of the dataset include" bullet (ilvl=0)
# ---------------------------------------------------------------------------
$rng1 = $d.Content
$found1 = $rng1.Find.Execute("Limitations include", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found1) { throw "Could not find 'Limitations include' paragraph" }
$para1 = $rng1.Paragraphs(1)
$target1 = $para1.Range
$target1.Collapse(0)
$xml1 = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>The majority of Kickstarter projects created that had highest success rate are in theater, music</w:t></w:r><w:r><w:t>, and</w:t></w:r><w:r><w:t xml:space="preserve"> film/video areas</w:t></w:r><w:r><w:t xml:space="preserve">. </w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Play projects were created most and with highest success outcome - 1066 plays projects out of a total of 4000 projects in all categories.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>US</w:t></w:r><w:r><w:t xml:space="preserve"> performed</w:t></w:r><w:r><w:t xml:space="preserve"> the highest numbers of Kickstarter projects</w:t></w:r><w:r><w:t>.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>The launch date is independent of the project outcome.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/></w:pPr><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Limitations</w:t></w:r><w:r><w:t xml:space="preserve"> of the dataset</w:t></w:r><w:r><w:t xml:space="preserve"> includ</w:t></w:r><w:r><w:t>e</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$target1.InsertXML($xml1)

# ---------------------------------------------------------------------------
# 2) Insert the 3 new "Limitations of the dataset" sub-bullets (+ trailing
#    blank ListParagraph) immediately before the "Other possible tables..."
#    bullet.
# ---------------------------------------------------------------------------
$rng2 = $d.Content
$found2 = $rng2.Find.Execute("Other possible tables and graphs that we could create include", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found2) { throw "Could not find 'Other possible tables...' paragraph" }
$para2 = $rng2.Paragraphs(1)
$target2 = $para2.Range
$target2.Collapse(1)
$xml2 = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>The dataset heavily depends on the outcomes of US as US had the highest number of projects done</w:t></w:r><w:r><w:t xml:space="preserve">. The dataset did not reveal trends for other </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t xml:space="preserve">countries </w:t></w:r><w:r><w:t>.</w:t></w:r><w:proofErr w:type="gramEnd"/></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">The dataset did not analyze the duration of each </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>project,</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> thus it was difficult to fully judge the success rate of each project.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>The average backer amount depends on the project' demographics. The dataset failed to address the demographic information of all projects within a country</w:t></w:r><w:r><w:t>.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/></w:pPr></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$target2.InsertXML($xml2)

# ---------------------------------------------------------------------------
# 3) The paragraph that used to follow "Other possible tables..." held the
#    original _GoBack bookmark; that bookmark now lives earlier in the
#    document (step 1), so this trailing paragraph is replaced with the new
#    "Graph the relationship..." sub-bullet, an indented blank paragraph, and
#    a final empty paragraph.
# ---------------------------------------------------------------------------
$rng3 = $d.Content
$found3 = $rng3.Find.Execute("Other possible tables and graphs that we could create include", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found3) { throw "Could not find 'Other possible tables...' paragraph (pass 2)" }
$para3 = $rng3.Paragraphs(1)
$trailingPara = $para3.Next()
$target3 = $trailingPara.Range
$target3.Collapse(0)
$xml3 = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Graph the relationship between the duration of the project and the success rate – longer project duration may increase the likelihood of meeting the project goal</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:ind w:left="1440"/></w:pPr></w:p><w:p/></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$target3.InsertXML($xml3)

# remove the (now-duplicate) old bookmark paragraph itself
$trailingPara.Range.Delete()

Write-Output "Edit complete. Paragraph count: $($d.Paragraphs.Count)"
